$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "A measure of how many standard deviations a value is to the right of the mean which is calculated by:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A measure of how many standard deviations a value is to the right of the mean which is calculated by:",
    2
)

$d.Content.Find.Execute(
    "Divide by n - 1 when an unbiased estimator of the popular variance is required.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Divide by n - 1 when an unbiased estimator of the popular variance is required.",
    2
)
